# Update the cryptos price/volume snapshot (GitHub Actions scheduled refresh).
# D column values that look purely numeric are prefixed with a leading
# apostrophe so Excel keeps them stored as text (matching the source feed's
# formatted strings, e.g. "0.999") instead of silently coercing them to
# floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.995.94'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '2.914.80'
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'589.83"
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').Value = "'144.94"
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').Value = "'6.93"
$ws.Range('E9').Value = '  +3.68%  '
$ws.Range('E10').Value = '  -2.14%  '
$ws.Range('E11').Value = '  -1.63%  '
$ws.Range('E12').Value = '  -0.49%  '
$ws.Range('D13').Value = "'33.38"
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('D15').Value = '3.393.53'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').Value = '60.879.62'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('D18').Value = '2.910.97'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').Value = "'432.99"
$ws.Range('E19').Value = '  +1.13%  '
$ws.Range('E20').Value = '  -1.38%  '
$ws.Range('E21').Value = '  -0.78%  '
$ws.Range('D22').Value = "'7.12"
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('D23').Value = "'81.41"
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('D24').Value = "'10.84"
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('D25').Value = "'2.21"
$ws.Range('E25').Value = '  -1.47%  '
$ws.Range('D26').Value = "'11.76"
$ws.Range('E26').Value = '  -1.10%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('E28').Value = '  +4.56%  '
$ws.Range('E29').Value = '  -0.82%  '
$ws.Range('D30').Value = "'6.95"
$ws.Range('E30').Value = '  -3.87%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').Value = "'0.109"
$ws.Range('E32').Value = '  +2.75%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').Value = '0.0₃0868'
$ws.Range('E34').Value = '  -0.81%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').Value = "'5.61"
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('D38').Value = "'1.97"
$ws.Range('E38').Value = '  -1.30%  '
$ws.Range('E39').Value = '  -4.12%  '
$ws.Range('E40').Value = '  -0.69%  '
$ws.Range('D41').Value = "'0.285"
$ws.Range('E41').Value = '  -3.71%  '
$ws.Range('D42').Value = "'41.17"
$ws.Range('E42').Value = '  -0.70%  '
$ws.Range('D43').Value = "'376.24"
$ws.Range('E43').Value = '  -0.43%  '
$ws.Range('E44').Value = '  -1.35%  '
$ws.Range('D45').Value = '2.691.36'
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').Value = "'133.25"
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('D48').Value = "'23.76"
$ws.Range('E48').Value = '  -2.64%  '
$ws.Range('D49').Value = "'0.105"
$ws.Range('E49').Value = '  -0.68%  '
$ws.Range('D50').Value = "'1.99"
$ws.Range('E50').Value = '  -2.41%  '
$ws.Range('E51').Value = '  -0.52%  '
